{"js": "// Load the first (only) table and all of its rows.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Load every row's cells so we can read/write each cell's text.\nfor (let i = 0; i < rows.items.length; i++) {\n  rows.items[i].cells.load(\"items\");\n}\nawait context.sync();\n\nfunction setRowText(rowIndex, text) {\n  rows.items[rowIndex].cells.items[0].value = text;\n}\n\n// Rows 0-2: \"99.87\" / \"0.14\" / \"106\" -> \"0M\" / \"0M\" / \"0M\"\nsetRowText(0, \"0M\");\nsetRowText(1, \"0M\");\nsetRowText(2, \"0M\");\n\n// Row 3: \"301\" -> \"700\"\nsetRowText(3, \"700\");\n\n// Row 5 (index 5): \"0.00032\" -> \"0.00026\"\nsetRowText(5, \"0.00026\");\n// Row 6: \"0.00014\" -> \"0.00034\"\nsetRowText(6, \"0.00034\");\n// Row 7: \"0.00004\" -> \"0.00040\"\nsetRowText(7, \"0.00040\");\n// Row 8: \"0.00015\" -> \"0.13711\" (rows 9, 10, 11 below it are removed)\nsetRowText(8, \"0.13711\");\n\n// Last 3 rows (43, 44, 45): collapse the tab-separated runs down to a\n// single short value each.\nconst lastIdx = rows.items.length - 1;\nsetRowText(lastIdx - 2, \"99.87\");\nsetRowText(lastIdx - 1, \"0.14\");\nsetRowText(lastIdx, \"106\");\n\nawait context.sync();\n\n// Insert 3 new single-value rows right after row 3 (\"700\").\nrows.items[3].insertRows(\"After\", 3, [[\"0.00002\"], [\"0.00065\"], [\"0.00017\"]]);\nawait context.sync();\n\n// Refresh the row collection (indices shifted after the insert above) and\n// delete the 3 rows that followed \"0.13711\" (originally rows 9, 10, 11;\n// now shifted down by 3 to indices 12, 13, 14).\nrows.load(\"items\");\nawait context.sync();\n\nrows.items[12].delete();\nrows.items[12].delete();\nrows.items[12].delete();\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# Rows 1-3 (1-based): \"99.87\" / \"0.14\" / \"106\" -> \"0M\" / \"0M\" / \"0M\"\n$table.Cell(1, 1).Range.Text = \"0M\"\n$table.Cell(2, 1).Range.Text = \"0M\"\n$table.Cell(3, 1).Range.Text = \"0M\"\n\n# Row 4: \"301\" -> \"700\"\n$table.Cell(4, 1).Range.Text = \"700\"\n\n# Insert 3 new single-value rows right after row 4 (\"700\"), i.e. before\n# the (then-current) row 5 (\"0.00004\"), in order so they read\n# 0.00002 / 0.00065 / 0.00017.\n$beforeRow = $table.Rows.Item(5)\n$newRow = $table.Rows.Add($beforeRow)\n$table.Cell($newRow.Index, 1).Range.Text = \"0.00002\"\n\n$beforeRow = $table.Rows.Item(6)\n$newRow = $table.Rows.Add($beforeRow)\n$table.Cell($newRow.Index, 1).Range.Text = \"0.00065\"\n\n$beforeRow = $table.Rows.Item(7)\n$newRow = $table.Rows.Add($beforeRow)\n$table.Cell($newRow.Index, 1).Range.Text = \"0.00017\"\n\n# Rows that were originally 6/7/8 (\"0.00032\"/\"0.00014\"/\"0.00004\") are now\n# shifted down by 3 to rows 9/10/11.\n$table.Cell(9, 1).Range.Text = \"0.00026\"\n$table.Cell(10, 1).Range.Text = \"0.00034\"\n$table.Cell(11, 1).Range.Text = \"0.00040\"\n\n# Row that was originally 9 (\"0.00015\") is now row 12; it becomes \"0.13711\"\n# and the 3 rows that used to follow it (originally 10/11/12, now 13/14/15)\n# are deleted entirely.\n$table.Cell(12, 1).Range.Text = \"0.13711\"\n$table.Rows.Item(13).Delete()\n$table.Rows.Item(13).Delete()\n$table.Rows.Item(13).Delete()\n\n# Last 3 rows: collapse the tab-separated runs down to a single short\n# value each.\n$rowCount = $table.Rows.Count\n$table.Cell($rowCount - 2, 1).Range.Text = \"99.87\"\n$table.Cell($rowCount - 1, 1).Range.Text = \"0.14\"\n$table.Cell($rowCount, 1).Range.Text = \"106\"\n"}
